# Update the "Budget" column (Q) so it reflects the new total budget of
# 80,000,000 being distributed according to each municipality's existing
# budgetShare (column P): Budget = budgetShare * 80,000,000.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTotalBudget = 80000000

for ($r = 2; $r -le 291; $r++) {
    $share = $ws.Cells.Item($r, 16).Value2
    $ws.Cells.Item($r, 17).Value2 = $share * $newTotalBudget
}
